$wb = $excel.ActiveWorkbook

# Reference to existing "Message" sheet, used as a template for the two
# new "*_Message" sheets.
$msg = $wb.Worksheets.Item("Message")

# ---------------------------------------------------------------------
# 1) CypherOutput_Message -- exact copy of the "Message" sheet (rows 1-10)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cypherMsg = $wb.Worksheets.Add($null, $lastSheet)
$cypherMsg.Name = "CypherOutput_Message"
for ($r = 1; $r -le 10; $r++) {
    $cypherMsg.Cells.Item($r, 1).Value = $msg.Cells.Item($r, 1).Value()
}

# ---------------------------------------------------------------------
# 2) StatOutput -- two-row stats table
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOut = $wb.Worksheets.Add($null, $lastSheet)
$statOut.Name = "StatOutput"
$statOut.Cells.Item(1, 1).Value = "number_of_files"
$statOut.Cells.Item(1, 2).Value = "number_of_sample"
$statOut.Cells.Item(1, 3).Value = "number_of_cases"
$statOut.Cells.Item(1, 4).Value = "number_of_study"
$statOut.Cells.Item(2, 1).Value = "'2"
$statOut.Cells.Item(2, 2).Value = "'5"
$statOut.Cells.Item(2, 3).Value = "'3"
$statOut.Cells.Item(2, 4).Value = "'2"

# ---------------------------------------------------------------------
# 3) StatOutput_Message -- "Message" sheet content repeated twice, with
#    the Cypher query text of the second block swapped for the
#    StatOutput query.
# ---------------------------------------------------------------------
$statQuery = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Other']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statMsg = $wb.Worksheets.Add($null, $lastSheet)
$statMsg.Name = "StatOutput_Message"
for ($r = 1; $r -le 10; $r++) {
    $statMsg.Cells.Item($r, 1).Value = $msg.Cells.Item($r, 1).Value()
}
for ($r = 1; $r -le 7; $r++) {
    $statMsg.Cells.Item($r + 10, 1).Value = $msg.Cells.Item($r, 1).Value()
}
$statMsg.Cells.Item(18, 1).Value = $statQuery
$statMsg.Cells.Item(19, 1).Value = $msg.Cells.Item(9, 1).Value()
$statMsg.Cells.Item(20, 1).Value = $msg.Cells.Item(10, 1).Value()
